$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows (2-30) down to (3-31)
$ws.Rows("2:2").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# The inserted row inherits the header row's formatting (bold/borders); reset it
# to match the plain formatting used by the rest of the data rows.
$ws.Rows("2:2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the newly inserted row 2 with the new record
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(2, 4).Value = 44691
$ws.Cells.Item(2, 5).Value = 15
$ws.Cells.Item(2, 6).Value = 100112045
$ws.Cells.Item(2, 7).Value = "Zapallo"
$ws.Cells.Item(2, 8).Value = "Camote"
$ws.Cells.Item(2, 9).Value = "1a (cosecha)"
$ws.Cells.Item(2, 10).Value = 700
$ws.Cells.Item(2, 11).Value = 580
$ws.Cells.Item(2, 12).Value = 600
$ws.Cells.Item(2, 13).Value = 590
$ws.Cells.Item(2, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(2, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(2, 16).Value = 590
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"
